$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 430 (shifts row 430 and below down by one)
$ws.Rows.Item(430).Insert()

# Populate the newly inserted row 430 with data
$ws.Cells.Item(430, 1).Value = 4
$ws.Cells.Item(430, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(430, 3).Value = 'Los Lagos'
$ws.Cells.Item(430, 4).Value = 45204
$ws.Cells.Item(430, 5).Value = 10
$ws.Cells.Item(430, 6).Value = 100112037
$ws.Cells.Item(430, 7).Value = 'Cebollín'
$ws.Cells.Item(430, 8).Value = 'Sin especificar'
$ws.Cells.Item(430, 9).Value = 'Primera'
$ws.Cells.Item(430, 10).Value = 90
$ws.Cells.Item(430, 11).Value = 6500
$ws.Cells.Item(430, 12).Value = 6500
$ws.Cells.Item(430, 13).Value = 6500
$ws.Cells.Item(430, 14).Value = '$/paquete 36 unidades'
$ws.Cells.Item(430, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(430, 16).Value = 181
$ws.Cells.Item(430, 17).Value = 36
$ws.Cells.Item(430, 18).Value = 'Hortaliza'
